# Scheduled-runner profit/price refresh across leve-crafting sheets.
# Columns: H=currentAveragePrice, I=currentAveragePriceNQ, J=currentAveragePriceHQ,
#          K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(20, 8).Value = 12378
$ws.Cells.Item(20, 9).Value = 633.3333
$ws.Cells.Item(20, 10).Value = 29995
$ws.Cells.Item(20, 11).Value = 633.3333
$ws.Cells.Item(20, 12).Value = 29995
$ws.Cells.Item(20, 13).Value = -403.3333
$ws.Cells.Item(20, 14).Value = -30455
$ws.Cells.Item(35, 8).Value = 12378
$ws.Cells.Item(35, 9).Value = 633.3333
$ws.Cells.Item(35, 10).Value = 29995
$ws.Cells.Item(35, 11).Value = 633.3333
$ws.Cells.Item(35, 12).Value = 29995
$ws.Cells.Item(35, 13).Value = -254.3333
$ws.Cells.Item(35, 14).Value = -30753
$ws.Cells.Item(40, 8).Value = 11948.8
$ws.Cells.Item(40, 10).Value = 2279
$ws.Cells.Item(40, 12).Value = 2279
$ws.Cells.Item(40, 14).Value = -2629
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 58825860
$ws.Cells.Item(45, 9).Value = 111113064
$ws.Cells.Item(45, 10).Value = 2758
$ws.Cells.Item(45, 11).Value = 111113064
$ws.Cells.Item(45, 12).Value = 2758
$ws.Cells.Item(45, 13).Value = -111112687
$ws.Cells.Item(45, 14).Value = -3512
$ws.Cells.Item(117, 8).Value = 60199.6
$ws.Cells.Item(117, 10).Value = 60199.6
$ws.Cells.Item(117, 12).Value = 60199.6
$ws.Cells.Item(117, 14).Value = -69377.60000000001
$ws.Cells.Item(119, 8).Value = 59899
$ws.Cells.Item(119, 10).Value = 59899
$ws.Cells.Item(119, 12).Value = 59899
$ws.Cells.Item(119, 14).Value = -69575
$ws.Cells.Item(120, 8).Value = 0
$ws.Cells.Item(120, 10).Value = 0
$ws.Cells.Item(120, 12).Value = 0
$ws.Cells.Item(120, 14).ClearContents()
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(117, 8).Value = 47719.184
$ws.Cells.Item(117, 10).Value = 47719.184
$ws.Cells.Item(117, 12).Value = 47719.184
$ws.Cells.Item(117, 14).Value = -56897.184
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(119, 14).ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(36, 8).Value = 3500
$ws.Cells.Item(36, 9).Value = 0
$ws.Cells.Item(36, 11).Value = 0
$ws.Cells.Item(36, 13).ClearContents()
$ws.Cells.Item(40, 8).Value = 3500
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(47, 8).Value = 23400
$ws.Cells.Item(47, 9).Value = 21250
$ws.Cells.Item(47, 10).Value = 32000
$ws.Cells.Item(47, 11).Value = 21250
$ws.Cells.Item(47, 12).Value = 32000
$ws.Cells.Item(47, 13).Value = -20684
$ws.Cells.Item(47, 14).Value = -33132
$ws.Cells.Item(99, 8).Value = 1930.1305
$ws.Cells.Item(99, 9).Value = 1917.4445
$ws.Cells.Item(99, 10).Value = 1938.2858
$ws.Cells.Item(99, 11).Value = 1917.4445
$ws.Cells.Item(99, 12).Value = 1938.2858
$ws.Cells.Item(99, 13).Value = -419.4445000000001
$ws.Cells.Item(99, 14).Value = -4934.2858
$ws.Cells.Item(100, 8).Value = 50000
$ws.Cells.Item(100, 10).Value = 50000
$ws.Cells.Item(100, 12).Value = 50000
$ws.Cells.Item(100, 14).Value = -52164
$ws.Cells.Item(119, 8).Value = 51999.5
$ws.Cells.Item(119, 10).Value = 51999.5
$ws.Cells.Item(119, 12).Value = 51999.5
$ws.Cells.Item(119, 14).Value = -61675.5
$ws.Cells.Item(121, 8).Value = 44008
$ws.Cells.Item(121, 10).Value = 44008
$ws.Cells.Item(121, 12).Value = 44008
$ws.Cells.Item(121, 14).Value = -46628
$ws.Cells.Item(122, 8).Value = 93028.84
$ws.Cells.Item(122, 9).Value = 93028.84
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 279086.52
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -276636.52
$ws.Cells.Item(122, 14).ClearContents()
$ws.Cells.Item(123, 8).Value = 25000
$ws.Cells.Item(123, 10).Value = 25000
$ws.Cells.Item(123, 12).Value = 25000
$ws.Cells.Item(123, 14).Value = -34800
$ws.Cells.Item(124, 8).Value = 40000
$ws.Cells.Item(124, 10).Value = 40000
$ws.Cells.Item(124, 12).Value = 40000
$ws.Cells.Item(124, 14).Value = -44910
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 1930.1305
$ws.Cells.Item(126, 9).Value = 1917.4445
$ws.Cells.Item(126, 10).Value = 1938.2858
$ws.Cells.Item(126, 11).Value = 5752.333500000001
$ws.Cells.Item(126, 12).Value = 5814.857400000001
$ws.Cells.Item(126, 13).Value = -3282.333500000001
$ws.Cells.Item(126, 14).Value = -10754.8574
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 1538.5
$ws.Cells.Item(46, 10).Value = 2977
$ws.Cells.Item(46, 12).Value = 8931
$ws.Cells.Item(46, 14).Value = -9113
$ws.Cells.Item(68, 8).Value = 11286197
$ws.Cells.Item(68, 9).Value = 10102171
$ws.Cells.Item(68, 10).Value = 11906401
$ws.Cells.Item(68, 11).Value = 30306513
$ws.Cells.Item(68, 12).Value = 35719203
$ws.Cells.Item(68, 13).Value = -30305702
$ws.Cells.Item(68, 14).Value = -35720825
$ws.Cells.Item(71, 8).Value = 11286197
$ws.Cells.Item(71, 9).Value = 10102171
$ws.Cells.Item(71, 10).Value = 11906401
$ws.Cells.Item(71, 11).Value = 90919539
$ws.Cells.Item(71, 12).Value = 107157609
$ws.Cells.Item(71, 13).Value = -90915483
$ws.Cells.Item(71, 14).Value = -107165721
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 14).ClearContents()
$ws.Cells.Item(70, 8).Value = 5940.909
$ws.Cells.Item(70, 9).Value = 5884.2104
$ws.Cells.Item(70, 10).Value = 6300
$ws.Cells.Item(70, 11).Value = 5884.2104
$ws.Cells.Item(70, 12).Value = 6300
$ws.Cells.Item(70, 13).Value = -5614.2104
$ws.Cells.Item(70, 14).Value = -6840
$ws.Cells.Item(73, 8).Value = 5940.909
$ws.Cells.Item(73, 9).Value = 5884.2104
$ws.Cells.Item(73, 10).Value = 6300
$ws.Cells.Item(73, 11).Value = 5884.2104
$ws.Cells.Item(73, 12).Value = 6300
$ws.Cells.Item(73, 13).Value = -4948.2104
$ws.Cells.Item(73, 14).Value = -8172
$ws.Cells.Item(113, 8).Value = 16668461
$ws.Cells.Item(113, 9).Value = 25001702
$ws.Cells.Item(113, 10).Value = 1981.5
$ws.Cells.Item(113, 11).Value = 25001702
$ws.Cells.Item(113, 12).Value = 1981.5
$ws.Cells.Item(113, 13).Value = -24999532
$ws.Cells.Item(113, 14).Value = -6321.5
$ws.Cells.Item(115, 8).Value = 0
$ws.Cells.Item(115, 10).Value = 0
$ws.Cells.Item(115, 12).Value = 0
$ws.Cells.Item(115, 14).ClearContents()
$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 25347.059
$ws.Cells.Item(126, 9).Value = 133637.33
$ws.Cells.Item(126, 10).Value = 2142
$ws.Cells.Item(126, 11).Value = 400911.99
$ws.Cells.Item(126, 12).Value = 6426
$ws.Cells.Item(126, 13).Value = -398441.99
$ws.Cells.Item(126, 14).Value = -11366
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(60, 8).Value = 30966.666
$ws.Cells.Item(60, 10).Value = 30966.666
$ws.Cells.Item(60, 12).Value = 30966.666
$ws.Cells.Item(60, 14).Value = -31984.666
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents()
$ws.Cells.Item(117, 8).Value = 0
$ws.Cells.Item(117, 10).Value = 0
$ws.Cells.Item(117, 12).Value = 0
$ws.Cells.Item(117, 14).ClearContents()
$ws.Cells.Item(118, 8).Value = 0
$ws.Cells.Item(118, 10).Value = 0
$ws.Cells.Item(118, 12).Value = 0
$ws.Cells.Item(118, 14).ClearContents()
$ws.Cells.Item(120, 8).Value = 55000
$ws.Cells.Item(120, 10).Value = 55000
$ws.Cells.Item(120, 12).Value = 55000
$ws.Cells.Item(120, 14).Value = -64676
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 14).ClearContents()
$ws.Cells.Item(47, 8).Value = 25000
$ws.Cells.Item(47, 10).Value = 25000
$ws.Cells.Item(47, 12).Value = 25000
$ws.Cells.Item(47, 14).Value = -26144
$ws.Cells.Item(126, 8).Value = 2452360.8
$ws.Cells.Item(126, 9).Value = 2452360.8
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 7357082.399999999
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -9447.5
$ws.Cells.Item(126, 14).ClearContents()
